$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 371530
$ws.Cells.Item(2, 4).Value = 478438275
$ws.Cells.Item(3, 3).Value = 296
$ws.Cells.Item(3, 4).Value = 358730
$ws.Cells.Item(4, 3).Value = 395
$ws.Cells.Item(4, 4).Value = 575460
$ws.Cells.Item(9, 3).Value = 985
$ws.Cells.Item(9, 4).Value = 1523938
$ws.Cells.Item(10, 3).Value = 21
$ws.Cells.Item(10, 4).Value = 48256
$ws.Cells.Item(11, 3).Value = 132060
$ws.Cells.Item(11, 4).Value = 202599608
$ws.Cells.Item(12, 3).Value = 155
$ws.Cells.Item(12, 4).Value = 230677
$ws.Cells.Item(13, 3).Value = 70005
$ws.Cells.Item(13, 4).Value = 107202458
$ws.Cells.Item(17, 3).Value = 4234
$ws.Cells.Item(17, 4).Value = 6065494
$ws.Cells.Item(22, 3).Value = 9284
$ws.Cells.Item(22, 4).Value = 13521155
$ws.Cells.Item(24, 3).Value = 87824
$ws.Cells.Item(24, 4).Value = 109419554
$ws.Cells.Item(25, 3).Value = 58
$ws.Cells.Item(25, 4).Value = 94209
$ws.Cells.Item(29, 3).Value = 307
$ws.Cells.Item(29, 4).Value = 440083
$ws.Cells.Item(30, 3).Value = 35854
$ws.Cells.Item(30, 4).Value = 53581734
$ws.Cells.Item(33, 3).Value = 13109
$ws.Cells.Item(33, 4).Value = 19521818
$ws.Cells.Item(36, 3).Value = 1663
$ws.Cells.Item(36, 4).Value = 2392578
$ws.Cells.Item(38, 3).Value = 2432
$ws.Cells.Item(38, 4).Value = 3490955
$ws.Cells.Item(39, 3).Value = 108987
$ws.Cells.Item(39, 4).Value = 136858251
$ws.Cells.Item(40, 3).Value = 70
$ws.Cells.Item(40, 4).Value = 82368
$ws.Cells.Item(41, 3).Value = 90
$ws.Cells.Item(41, 4).Value = 126078
$ws.Cells.Item(45, 3).Value = 951
$ws.Cells.Item(45, 4).Value = 1404517
$ws.Cells.Item(47, 3).Value = 48847
$ws.Cells.Item(47, 4).Value = 72731871
$ws.Cells.Item(49, 3).Value = 10455
$ws.Cells.Item(49, 4).Value = 15456034
$ws.Cells.Item(51, 3).Value = 1512
$ws.Cells.Item(51, 4).Value = 2111195
$ws.Cells.Item(54, 3).Value = 3078
$ws.Cells.Item(54, 4).Value = 4436798
$ws.Cells.Item(55, 3).Value = 78465
$ws.Cells.Item(55, 4).Value = 99007845
$ws.Cells.Item(57, 3).Value = 57
$ws.Cells.Item(57, 4).Value = 89637
$ws.Cells.Item(60, 3).Value = 406
$ws.Cells.Item(60, 4).Value = 595015
$ws.Cells.Item(62, 3).Value = 31588
$ws.Cells.Item(62, 4).Value = 47359833
$ws.Cells.Item(63, 3).Value = 32
$ws.Cells.Item(63, 4).Value = 44450
$ws.Cells.Item(65, 3).Value = 12882
$ws.Cells.Item(65, 4).Value = 19145030
$ws.Cells.Item(67, 3).Value = 1469
$ws.Cells.Item(67, 4).Value = 2058169
$ws.Cells.Item(71, 3).Value = 1982
$ws.Cells.Item(71, 4).Value = 2909419
$ws.Cells.Item(73, 3).Value = 23648
$ws.Cells.Item(73, 4).Value = 31163970
$ws.Cells.Item(77, 3).Value = 8660
$ws.Cells.Item(77, 4).Value = 13164204
$ws.Cells.Item(79, 3).Value = 5971
$ws.Cells.Item(79, 4).Value = 9022672
$ws.Cells.Item(80, 3).Value = 582
$ws.Cells.Item(80, 4).Value = 839491
$ws.Cells.Item(81, 3).Value = 388
$ws.Cells.Item(81, 4).Value = 572028
$ws.Cells.Item(82, 3).Value = 160494
$ws.Cells.Item(82, 4).Value = 200523980
$ws.Cells.Item(84, 3).Value = 98
$ws.Cells.Item(84, 4).Value = 143578
$ws.Cells.Item(86, 3).Value = 480
$ws.Cells.Item(86, 4).Value = 709352
$ws.Cells.Item(88, 3).Value = 70351
$ws.Cells.Item(88, 4).Value = 104907477
$ws.Cells.Item(91, 3).Value = 33691
$ws.Cells.Item(91, 4).Value = 50157610
$ws.Cells.Item(93, 3).Value = 3060
$ws.Cells.Item(93, 4).Value = 4477236
$ws.Cells.Item(95, 3).Value = 3747
$ws.Cells.Item(95, 4).Value = 5466073
$ws.Cells.Item(96, 3).Value = 39356
$ws.Cells.Item(96, 4).Value = 54054988
$ws.Cells.Item(100, 3).Value = 9799
$ws.Cells.Item(100, 4).Value = 14925015
$ws.Cells.Item(102, 3).Value = 9098
$ws.Cells.Item(102, 4).Value = 13685471
$ws.Cells.Item(104, 3).Value = 613
$ws.Cells.Item(104, 4).Value = 878127
$ws.Cells.Item(106, 3).Value = 17206
$ws.Cells.Item(106, 4).Value = 32992757
$ws.Cells.Item(109, 3).Value = 3956
$ws.Cells.Item(109, 4).Value = 8140751
$ws.Cells.Item(111, 3).Value = 5658
$ws.Cells.Item(111, 4).Value = 11823621
$ws.Cells.Item(114, 3).Value = 336
$ws.Cells.Item(114, 4).Value = 675578
$ws.Cells.Item(116, 3).Value = 163083
$ws.Cells.Item(116, 4).Value = 203707585
$ws.Cells.Item(120, 3).Value = 1041
$ws.Cells.Item(120, 4).Value = 1578470
$ws.Cells.Item(122, 3).Value = 59463
$ws.Cells.Item(122, 4).Value = 90839555
$ws.Cells.Item(123, 3).Value = 111
$ws.Cells.Item(123, 4).Value = 171171
$ws.Cells.Item(124, 3).Value = 32378
$ws.Cells.Item(124, 4).Value = 49490443
$ws.Cells.Item(125, 3).Value = 1406
$ws.Cells.Item(125, 4).Value = 1961480
$ws.Cells.Item(129, 3).Value = 3073
$ws.Cells.Item(129, 4).Value = 4526822
$ws.Cells.Item(131, 3).Value = 672564
$ws.Cells.Item(131, 4).Value = 924790001
$ws.Cells.Item(132, 3).Value = 114
$ws.Cells.Item(132, 4).Value = 167262
$ws.Cells.Item(133, 3).Value = 269
$ws.Cells.Item(133, 4).Value = 447901
$ws.Cells.Item(136, 3).Value = 1753
$ws.Cells.Item(136, 4).Value = 2940143
$ws.Cells.Item(137, 3).Value = 40
$ws.Cells.Item(137, 4).Value = 79510
$ws.Cells.Item(138, 3).Value = 253040
$ws.Cells.Item(138, 4).Value = 401749785
$ws.Cells.Item(139, 3).Value = 601
$ws.Cells.Item(139, 4).Value = 1182549
$ws.Cells.Item(140, 3).Value = 24
$ws.Cells.Item(140, 4).Value = 53957
$ws.Cells.Item(141, 3).Value = 241101
$ws.Cells.Item(141, 4).Value = 385332368
$ws.Cells.Item(144, 3).Value = 3153
$ws.Cells.Item(144, 4).Value = 4548870
$ws.Cells.Item(147, 3).Value = 9273
$ws.Cells.Item(147, 4).Value = 13863695
$ws.Cells.Item(150, 3).Value = 51456
$ws.Cells.Item(150, 4).Value = 69206237
$ws.Cells.Item(154, 3).Value = 9
$ws.Cells.Item(154, 4).Value = 19118
$ws.Cells.Item(156, 3).Value = 15731
$ws.Cells.Item(156, 4).Value = 23403390
$ws.Cells.Item(157, 3).Value = 4295
$ws.Cells.Item(157, 4).Value = 6281513
$ws.Cells.Item(162, 3).Value = 538
$ws.Cells.Item(162, 4).Value = 791921
$ws.Cells.Item(163, 3).Value = 20664
$ws.Cells.Item(163, 4).Value = 27545267
$ws.Cells.Item(167, 3).Value = 8588
$ws.Cells.Item(167, 4).Value = 12914262
$ws.Cells.Item(169, 3).Value = 6122
$ws.Cells.Item(169, 4).Value = 9136021
$ws.Cells.Item(172, 3).Value = 349
$ws.Cells.Item(172, 4).Value = 519444
$ws.Cells.Item(174, 3).Value = 29719
$ws.Cells.Item(174, 4).Value = 61277131
$ws.Cells.Item(176, 3).Value = 3000
$ws.Cells.Item(176, 4).Value = 6200380
$ws.Cells.Item(177, 3).Value = 402
$ws.Cells.Item(177, 4).Value = 825982
$ws.Cells.Item(179, 3).Value = 94
$ws.Cells.Item(179, 4).Value = 199530
$ws.Cells.Item(180, 3).Value = 189
$ws.Cells.Item(180, 4).Value = 415887
$ws.Cells.Item(181, 3).Value = 99350
$ws.Cells.Item(181, 4).Value = 124412865
$ws.Cells.Item(186, 3).Value = 673
$ws.Cells.Item(186, 4).Value = 1017500
$ws.Cells.Item(188, 3).Value = 37876
$ws.Cells.Item(188, 4).Value = 57267264
$ws.Cells.Item(190, 3).Value = 15252
$ws.Cells.Item(190, 4).Value = 23127385
$ws.Cells.Item(192, 3).Value = 1340
$ws.Cells.Item(192, 4).Value = 1887054
$ws.Cells.Item(194, 3).Value = 2188
$ws.Cells.Item(194, 4).Value = 3183099
$ws.Cells.Item(196, 3).Value = 269725
$ws.Cells.Item(196, 4).Value = 335896009
$ws.Cells.Item(198, 3).Value = 196
$ws.Cells.Item(198, 4).Value = 281155
$ws.Cells.Item(202, 3).Value = 977
$ws.Cells.Item(202, 4).Value = 1470026
$ws.Cells.Item(204, 3).Value = 96156
$ws.Cells.Item(204, 4).Value = 143794046
$ws.Cells.Item(207, 3).Value = 38064
$ws.Cells.Item(207, 4).Value = 56437565
$ws.Cells.Item(210, 3).Value = 5582
$ws.Cells.Item(210, 4).Value = 8047721
$ws.Cells.Item(213, 3).Value = 6423
$ws.Cells.Item(213, 4).Value = 9090951
$ws.Cells.Item(216, 3).Value = 302985
$ws.Cells.Item(216, 4).Value = 378000974
$ws.Cells.Item(223, 3).Value = 683
$ws.Cells.Item(223, 4).Value = 1057023
$ws.Cells.Item(225, 3).Value = 107201
$ws.Cells.Item(225, 4).Value = 163869518
$ws.Cells.Item(228, 3).Value = 59916
$ws.Cells.Item(228, 4).Value = 90949260
$ws.Cells.Item(231, 3).Value = 5042
$ws.Cells.Item(231, 4).Value = 7124092
$ws.Cells.Item(234, 3).Value = 8099
$ws.Cells.Item(234, 4).Value = 11638662
$ws.Cells.Item(237, 3).Value = 121145
$ws.Cells.Item(237, 4).Value = 151574078
$ws.Cells.Item(239, 3).Value = 94
$ws.Cells.Item(239, 4).Value = 131901
$ws.Cells.Item(240, 3).Value = 14
$ws.Cells.Item(240, 4).Value = 18647
$ws.Cells.Item(242, 3).Value = 614
$ws.Cells.Item(242, 4).Value = 888837
$ws.Cells.Item(244, 3).Value = 55050
$ws.Cells.Item(244, 4).Value = 81774514
$ws.Cells.Item(246, 3).Value = 14683
$ws.Cells.Item(246, 4).Value = 21777042
$ws.Cells.Item(248, 3).Value = 1997
$ws.Cells.Item(248, 4).Value = 2903697
$ws.Cells.Item(250, 3).Value = 3263
$ws.Cells.Item(250, 4).Value = 4672030
$ws.Cells.Item(251, 3).Value = 301566
$ws.Cells.Item(251, 4).Value = 386977432
$ws.Cells.Item(252, 3).Value = 209
$ws.Cells.Item(252, 4).Value = 266658
$ws.Cells.Item(253, 3).Value = 283
$ws.Cells.Item(253, 4).Value = 402541
$ws.Cells.Item(258, 3).Value = 982
$ws.Cells.Item(258, 4).Value = 1524245
$ws.Cells.Item(260, 3).Value = 110620
$ws.Cells.Item(260, 4).Value = 171715362
$ws.Cells.Item(263, 3).Value = 78981
$ws.Cells.Item(263, 4).Value = 122830034
$ws.Cells.Item(265, 3).Value = 2662
$ws.Cells.Item(265, 4).Value = 3832416
$ws.Cells.Item(268, 3).Value = 6535
$ws.Cells.Item(268, 4).Value = 9640490
